# SSDM-12286 Fixed letter case inconsistencies.
# The "Vocabulary Code" header (column H, row 4) should read "Vocabulary code".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H4").Value = "Vocabulary code"

# Reflect the new active selection (H4) as seen in the target workbook.
$ws.Range("H4").Select()
